$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Section header rows: make column E match the shaded style already
# used by the other columns (C/D/F/G) in the same row -----------------
$ws.Range("C7").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$ws.Range("C17").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("C26").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("C30").Copy()
$ws.Range("E30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Assign "THÀNH VIÊN" (member) names for the task rows ------------
$ws.Range("E8").Value = "Lê Nguyễn Hoài Đăng"

$ws.Range("E18").Value = "Lê Phước Anh Đạt"
$ws.Range("E19").Value = "Lê Phước Anh Đạt"
$ws.Range("E20").Value = "Lê Phước Anh Đạt"
$ws.Range("E21").Value = "Lê Phước Anh Đạt"
$ws.Range("E22").Value = "Lê Phước Anh Đạt"

$ws.Range("E23").Value = "Cả nhóm"
$ws.Range("E24").Value = "Cả nhóm"
$ws.Range("E25").Value = "Cả nhóm"

# --- Update the view: scroll position + current selection ------------
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("E31").Select()
